# Daily attendance processing - normalize "Recorded By" (column G) ordering.
#
# For every data row, the "Recorded By" cell holds a comma-separated list of
# recorder identities (e.g. "System, dnasr281@gmail.com"). This pass rotates
# the list one position to the left (the first-listed recorder is moved to
# the end) - EXCEPT when the list already starts with dnasr281@gmail.com,
# which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $parts = $raw -split ", "

    if ($parts.Length -le 1) {
        continue
    }

    $first = $parts[0]

    if ($first.ToLower() -eq "dnasr281@gmail.com") {
        continue
    }

    $rest = $parts[1..($parts.Length - 1)]
    $rotated = $rest + @($first)
    $newVal = [string]::Join(", ", $rotated)

    $cell.Value = $newVal
}
